$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Trim the stray leading space from these Kannada (kan) and Hindi (hin) location-name cells.
$ws.Range("C122").Value = "ನನ್ನ ದೇಶ"
$ws.Range("C128").Value = "ಮ್ನಸ್ರ"
$ws.Range("C134").Value = "ಮೆಹದಿಯಾ"
$ws.Range("C136").Value = "ಎಣ್ಣೆಯ ಔಜಿ"
$ws.Range("C138").Value = "ಸಿಡಿ ತೈಬಿ"
$ws.Range("C140").Value = "ಸಿಡಿ ಅಲ್ಲಲ್ ತಾಜಿ"
$ws.Range("C142").Value = "ರಬತ್"
$ws.Range("C143").Value = "ರಬತ್"
$ws.Range("C152").Value = "ಹೇ ರಿಯಾದ್"
$ws.Range("C162").Value = "मेरा देश"
$ws.Range("C163").Value = "रबात बिक्री केनित्र"
$ws.Range("C164").Value = "केनिट्रा"
$ws.Range("C165").Value = "केनिट्रा"
$ws.Range("C166").Value = "बेन मंसूर"
$ws.Range("C170").Value = "माइग्रेन"
$ws.Range("C174").Value = "मेहदिया"
$ws.Range("C188").Value = "सौइसी"
$ws.Range("C192").Value = "हे रियादो"

# Restore the sheet view/selection state recorded in the workbook (scroll position + active cell).
$ws.Range("C122").Select()
$excel.ActiveWindow.ScrollRow = 118
$excel.ActiveWindow.ScrollColumn = 1
